$d = $word.ActiveDocument
$d.Content.Find.Execute("svn.  John Paul", $true, $false, $false, $false, $false,
                         $true, 1, $false, "svn client over Google code set up by William Peckham.  John Paul", 2)
